$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'256.62"
$ws.Range("E2").Formula = "'-1.27%"
$ws.Range("D3").Formula = "'27.27"
$ws.Range("E3").Formula = "'-2.58%"
$ws.Range("D4").Formula = "'4.557"
$ws.Range("E4").Formula = "'-12.68%"
$ws.Range("D5").Formula = "'0.05901"
$ws.Range("E5").Formula = "'-0.42%"
$ws.Range("E6").Formula = "'-1.51%"
$ws.Range("D7").Formula = "'0.8591"
$ws.Range("E7").Formula = "'-1.68%"
$ws.Range("D8").Formula = "'0.9326"
$ws.Range("E8").Formula = "'-7.64%"
$ws.Range("E9").Formula = "'-1.10%"
$ws.Range("D10").Formula = "'0.03650"
$ws.Range("E10").Formula = "'0.36%"
$ws.Range("E11").Formula = "'-1.98%"
$ws.Range("D12").Formula = "'0.03233"
$ws.Range("E12").Formula = "'1.12%"
$ws.Range("D13").Formula = "'0.09208"
$ws.Range("E13").Formula = "'-0.38%"
$ws.Range("D14").Formula = "'0.001545"
$ws.Range("E14").Formula = "'0.44%"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Formula = "'0.006104"
$ws.Range("E15").Formula = "'3.99%"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Formula = "'3.516"
$ws.Range("E16").Formula = "'0.53%"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Formula = "'3.190"
$ws.Range("E17").Formula = "'-1.32%"
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Formula = "'2.202"
$ws.Range("E18").Formula = "'-0.23%"
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").Formula = "'0.01039"
$ws.Range("E19").Formula = "'1,621.06%"
$ws.Range("E20").Formula = "'-2.15%"
$ws.Range("E21").Formula = "'-0.92%"
$ws.Range("D22").Formula = "'3.847"
$ws.Range("E22").Formula = "'9.43%"
$ws.Range("D23").Formula = "'0.04206"
$ws.Range("E23").Formula = "'0.65%"
$ws.Range("D24").Formula = "'0.001222"
$ws.Range("E24").Formula = "'0.49%"
$ws.Range("D25").Formula = "'0.004276"
$ws.Range("E25").Formula = "'-6.53%"
$ws.Range("E26").Formula = "'0.16%"
$ws.Range("D27").Formula = "'0.0001511"
$ws.Range("E27").Formula = "'-21.87%"
$ws.Range("D40").Formula = "'0.03829"
$ws.Range("E40").Formula = "'-0.42%"
$ws.Range("D41").Formula = "'0.006232"
$ws.Range("E41").Formula = "'55.44%"
$ws.Range("D42").Formula = "'0.1099"
$ws.Range("E42").Formula = "'-1.00%"
$ws.Range("E43").Formula = "'-7.25%"
$ws.Range("D44").Formula = "'0.01134"
$ws.Range("E44").Formula = "'3.93%"
$ws.Range("D45").Formula = "'0.00005453"
$ws.Range("E45").Formula = "'1.00%"
$ws.Range("E46").Formula = "'0.26%"
$ws.Range("D47").Formula = "'0.08802"
$ws.Range("E47").Formula = "'3.18%"
$ws.Range("D48").Formula = "'0.09732"
$ws.Range("E48").Formula = "'4,456.19%"
$ws.Range("E49").Formula = "'0.26%"
$ws.Range("E50").Formula = "'0.26%"
